$d = $word.ActiveDocument

$replacements = @(
    @("805×8=6440", "153×4=612"),
    @("685×6=4110", "154×3=462"),
    @("656×3=1968", "946×8=7568"),
    @("784×3=2352", "417×9=3753"),
    @("296×8=2368", "858×6=5148"),
    @("434×8=3472", "567×6=3402"),
    @("630×9=5670", "410×6=2460"),
    @("195×9=1755", "733×8=5864"),
    @("252×7=1764", "504×2=1008"),
    @("706×7=4942", "476×3=1428"),
    @("153×5=765",  "401×7=2807"),
    @("701×7=4907", "349×5=1745"),
    @("832×8=6656", "862×4=3448"),
    @("763×7=5341", "209×7=1463"),
    @("632×3=1896", "878×3=2634"),
    @("770×8=6160", "723×5=3615"),
    @("412×2=824",  "157×9=1413"),
    @("657×7=4599", "962×7=6734"),
    @("132×5=660",  "556×4=2224"),
    @("945×7=6615", "617×4=2468"),
    @("770×9=6930", "492×5=2460"),
    @("943×5=4715", "174×6=1044"),
    @("937×7=6559", "143×5=715"),
    @("513×4=2052", "607×8=4856"),
    @("980×8=7840", "145×4=580")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
